$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift header text D->E->F->G, put new label in D ---
# Capture old header text before overwriting
$oldD1 = $ws.Range("D1").Value()
$oldE1 = $ws.Range("E1").Value()
$oldF1 = $ws.Range("F1").Value()

$ws.Range("G1").Value = $oldF1
$ws.Range("F1").Value = $oldE1
$ws.Range("E1").Value = $oldD1
$ws.Range("D1").Value = "StartTime matlab datenum"

# --- Data rows 2-18: shift D,E,F -> E,F,G; blank out D (format matches col C) ---
for ($r = 2; $r -le 18; $r++) {
    $dVal = $ws.Range("D$r").Value()
    $eVal = $ws.Range("E$r").Value()
    $fVal = $ws.Range("F$r").Value()

    $ws.Range("G$r").Value = $fVal
    $ws.Range("F$r").Value = $eVal
    $ws.Range("E$r").Value = $dVal

    # Clear D and give it the same number format as column C on this row
    $ws.Range("C$r").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("D$r").Value = $null
}

$excel.CutCopyMode = 0

# --- New column H (rows 11-18): deviation from row 11 baseline ---
for ($r = 11; $r -le 18; $r++) {
    $ws.Range("H$r").Formula = '=G' + $r + '-$G$11'
}

# --- Column width: column G should match D:F width (19.77) ---
$ws.Columns("G").ColumnWidth = 19

# --- Page margins (points; 1.025in = 73.8pt) ---
$ws.PageSetup.TopMargin = 73.8
$ws.PageSetup.BottomMargin = 73.8

# --- Header / footer text ---
$ws.PageSetup.CenterHeader = "&A"
$ws.PageSetup.CenterFooter = "Page &P"

# --- Selection ---
$ws.Range("H12").Select()
